$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 90
$ws1.Range("F5").Value = 36
$ws1.Range("F8").Value = 8130
$ws1.Range("F9").Value = 762
$ws1.Range("F10").Value = 250
$ws1.Range("F11").Value = 1106
$ws1.Range("F12").Value = 798
$ws1.Range("F13").Value = 40
$ws1.Range("F14").Value = 33
$ws1.Range("F16").Value = 74
$ws1.Range("F19").Value = 868

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 90
$ws4.Range("F5").Value = 36
$ws4.Range("F9").Value = 8130
$ws4.Range("F10").Value = 762
$ws4.Range("F11").Value = 250
$ws4.Range("F12").Value = 1106
$ws4.Range("F13").Value = 798
$ws4.Range("F14").Value = 40
$ws4.Range("F15").Value = 33
$ws4.Range("F17").Value = 74
$ws4.Range("F20").Value = 868
